# Schedule - Spring 2020: shift the "date" column (B) up by one class session
# (each row now shows the date that used to belong to the following row),
# and the last row (27) loses its date since there is no further row to pull
# from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the class-date values in column B (rows 7-26) ----------------
$ws.Range("B7").Value  = 43873
$ws.Range("B8").Value  = 43878
$ws.Range("B9").Value  = 43880
$ws.Range("B10").Value = 43885
$ws.Range("B11").Value = 43887
$ws.Range("B12").Value = 43892
$ws.Range("B13").Value = 43894
$ws.Range("B14").Value = 43899
$ws.Range("B15").Value = 43901
$ws.Range("B16").Value = 43913
$ws.Range("B17").Value = 43915
$ws.Range("B18").Value = 43920
$ws.Range("B19").Value = 43922
$ws.Range("B20").Value = 43927
$ws.Range("B21").Value = 43929
$ws.Range("B22").Value = 43934
$ws.Range("B23").Value = 43936
$ws.Range("B24").Value = 43941
$ws.Range("B25").Value = 43943
$ws.Range("B26").Value = 43948

# Row 27 no longer carries a date - remove the cell entirely (not just its
# contents) so the row only keeps its C/D/E/F values, matching the source
# file where the trailing <c r="B27"> node is gone.
$ws.Range("B27").Clear()

# --- Update the view state (scroll position / active selection) ----------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("B27").Select()
